$wb = $excel.ActiveWorkbook

# The edits apply to the "tough_levels" sheet (second tab), which is also
# the tab that is active/selected in the workbook.
$ws = $wb.Worksheets.Item("tough_levels")
$ws.Activate()

# Buff the "strength" column (D) values for the two enemy level rows.
$ws.Range("D2").Value = 10
$ws.Range("D3").Value = 30

# Move the active selection as recorded in the saved file.
$ws.Range("D11").Select()
